# Add the U238 decay-chain sheet (mirrors Th232 / U235 layout) after U235.

$wb = $excel.ActiveWorkbook
$u235 = $wb.Worksheets.Item("U235")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "U238"

# ---------------------------------------------------------------------
# 1) Formats first (reuse U235's existing style indices; PasteSpecial of
#    formats-only keeps the shared cellXfs entries instead of minting new
#    ones).
# ---------------------------------------------------------------------

# Header row (row 1) - green header style
$u235.Range("A1:G1").Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)

# Data block rows 2-19, columns C:F always use the scientific-format style
$u235.Range("C2:F2").Copy()
$ws.Range("C2:F19").PasteSpecial(-4122)

# Column B uses the scientific-format style only on rows 2,4,5,6,7
$u235.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B4:B7").PasteSpecial(-4122)

# Column G only has data in row 2
$u235.Range("G2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

# Second header row (row 24) - orange header style
$u235.Range("A24:C24").Copy()
$ws.Range("A24:C24").PasteSpecial(-4122)

# Branching-ratio block rows 25-28; C26 carries the scientific-format style
$u235.Range("C26").Copy()
$ws.Range("C26").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Header text (reuses the existing shared-string entries)
# ---------------------------------------------------------------------

$ws.Range("A1").Value = "Decay"
$ws.Range("B1").Value = "Half-Life (yrs)"
$ws.Range("C1").Value = "Mean Life (yrs)"
$ws.Range("D1").Value = "Lambda"
$ws.Range("E1").Value = "Mean Life (Scaled)"
$ws.Range("F1").Value = "Lambda (Scaled)"
$ws.Range("G1").Value = "Lambda1"

$ws.Range("A24").Value = "Decay"
$ws.Range("B24").Value = "branching ratio beta (%)"
$ws.Range("C24").Value = "branching ratio alpha (%)"

# ---------------------------------------------------------------------
# 3) Column A - decay step indices (rows 2-19), literal integers
# ---------------------------------------------------------------------

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18

# ---------------------------------------------------------------------
# 4) Column B - half-lives (mix of literal values and formulas)
# ---------------------------------------------------------------------

$ws.Range("B2").Value = 4468000000
$ws.Range("B3").Formula = "=24.1/365"
$ws.Range("B4").Formula = "=1.159/60/24/365"
$ws.Range("B5").Value = 245500
$ws.Range("B6").Value = 75400
$ws.Range("B7").Value = 1600
$ws.Range("B8").Formula = "=3.8235/365"
$ws.Range("B9").Formula = "=3.098/60/24/365"
$ws.Range("B10").Formula = "=1.5/60/60/24/365"
$ws.Range("B11").Formula = "=27.06/60/24/365"
$ws.Range("B12").Formula = "=19.9/60/24/365"
$ws.Range("B13").Formula = "=0.0001636/60/60/24/365"
$ws.Range("B14").Formula = "=1.3/60/24/365"
$ws.Range("B15").Value = 22.2
$ws.Range("B16").Formula = "=5.012/365"
$ws.Range("B17").Formula = "=138.376/365"
$ws.Range("B18").Formula = "=8.32/60/24/365"
$ws.Range("B19").Formula = "=4.202/60/24/365"

# ---------------------------------------------------------------------
# 5) Column G - scaling constant (row 2 only)
# ---------------------------------------------------------------------

$ws.Range("G2").Value = 0.00000000015500000000000001

# ---------------------------------------------------------------------
# 6) Columns C:F - derived formulas
#    Row 2 is unique (absolute $G$2 anchor); rows 3-19 share one pattern
#    each (Excel will recreate the shared-formula group automatically).
# ---------------------------------------------------------------------

$ws.Range("C2").Formula = "=B2/LN(2)"
$ws.Range("D2").Formula = "=1/C2"
$ws.Range("F2").Formula = "=D2/`$G`$2"
$ws.Range("E2").Formula = "=1/F2"

$ws.Range("C3:C19").Formula = "=B3/LN(2)"
$ws.Range("D3:D19").Formula = "=1/C3"
$ws.Range("F3:F19").Formula = "=D3/`$G`$2"
$ws.Range("E3:E19").Formula = "=1/F3"

# ---------------------------------------------------------------------
# 7) Branching-ratio block (rows 25-28)
# ---------------------------------------------------------------------

$ws.Range("A25").Value = 8
$ws.Range("B25").Formula = "=0.0002*100"
$ws.Range("C25").Formula = "=0.9998*100"

$ws.Range("A26").Value = 11
$ws.Range("B26").Formula = "=0.9998*100"
$ws.Range("C26").Formula = "=0.0002*100"

$ws.Range("A27").Value = 14
$ws.Range("B27").Formula = "=1*100"
$ws.Range("C27").Formula = "=0.000000019*100"

$ws.Range("A28").Value = 15
$ws.Range("B28").Formula = "=1*100"
$ws.Range("C28").Formula = "=0.0000013*100"

# ---------------------------------------------------------------------
# 8) Sheet view niceties matching the authored file
# ---------------------------------------------------------------------

$ws.Range("C29").Select()
$excel.ActiveWindow.Zoom = 117
